$d = $word.ActiveDocument

# --- Step 1: paragraph 1 keeps its text, but the paragraph-mark run
#     properties (w:pPr/w:rPr) drop <w:rFonts w:hint="cs"/> and <w:rtl/>,
#     leaving only <w:lang w:bidi="fa-IR"/>. The run's own w:rPr is
#     untouched. Rebuild the paragraph via InsertXML targeted at its own
#     Range so nothing outside paragraph 1 is disturbed, keeping the
#     original rsid attributes on <w:p>.
$p1 = $d.Paragraphs.Item(1)
$para1Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
    '<w:p w:rsidR="00E224ED" w:rsidRDefault="00FA211B" w:rsidP="00FA211B">' +
      '<w:pPr>' +
        '<w:bidi/>' +
        '<w:rPr>' +
          '<w:lang w:bidi="fa-IR"/>' +
        '</w:rPr>' +
      '</w:pPr>' +
      '<w:r>' +
        '<w:rPr>' +
          '<w:rFonts w:hint="cs"/>' +
          '<w:rtl/>' +
          '<w:lang w:bidi="fa-IR"/>' +
        '</w:rPr>' +
        '<w:t>این پروژه توسط گروهی زبده طراحی و تست شد.</w:t>' +
      '</w:r>' +
    '</w:p>' +
  '</w:body>' +
  '</w:wordDocument>'
$p1.Range.InsertXML($para1Xml)

# --- Step 2: add a brand-new paragraph right after paragraph 1 holding
#     "asdfalfsjsfjjn", and relocate the _GoBack bookmark onto the end
#     of this new paragraph (it previously sat at the end of paragraph
#     1). Split first with InsertParagraphAfter so paragraph 1 is left
#     completely intact (and without the bookmark), then fill the new,
#     empty paragraph via InsertXML.
$p1 = $d.Paragraphs.Item(1)
$splitPoint = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$splitPoint.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(2)
$para2Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
    '<w:p>' +
      '<w:pPr>' +
        '<w:bidi/>' +
        '<w:rPr>' +
          '<w:rtl/>' +
          '<w:lang w:bidi="fa-IR"/>' +
        '</w:rPr>' +
      '</w:pPr>' +
      '<w:r>' +
        '<w:rPr>' +
          '<w:lang w:bidi="fa-IR"/>' +
        '</w:rPr>' +
        '<w:t>asdfalfsjsfjjn</w:t>' +
      '</w:r>' +
      '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
      '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>' +
  '</w:body>' +
  '</w:wordDocument>'
$p2.Range.InsertXML($para2Xml)
